# Update column G ("K" - strike count) values on Sheet1 rows 2-36
# per regenerated save_data (K instead of Strike#).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$gValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 3
    6  = 0
    7  = 2
    8  = 0
    9  = 3
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 3
    15 = 0
    16 = 0
    17 = 3
    18 = 1
    19 = 0
    20 = 2
    21 = 0
    22 = 2
    23 = 1
    24 = 3
    25 = 1
    26 = 2
    27 = 2
    28 = 5
    29 = 0
    30 = 2
    31 = 2
    32 = 3
    33 = 2
    34 = 2
    35 = 1
    36 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
